$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (dotted thousands separators, e.g. "95.246.19").
# These look numeric to Excel's auto-detection, so they are entered with a
# leading apostrophe to force text, matching the original inlineStr cell type.

# Row 2
$ws.Range("D2").Value = "'95.246.19"
$ws.Range("E2").Value = "  +2.12%  "

# Row 3
$ws.Range("D3").Value = "'3.618.01"
$ws.Range("E3").Value = "  +6.09%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("E5").Value = "  +2.09%  "

# Row 6
$ws.Range("D6").Value = "'656.33"
$ws.Range("E6").Value = "  +6.03%  "

# Row 7
$ws.Range("E7").Value = "  +1.93%  "

# Row 8
$ws.Range("D8").Value = "'0.405"
$ws.Range("E8").Value = "  +3.70%  "

# Row 9
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.04%  "

# Row 10
$ws.Range("D10").Value = "'0.992"
$ws.Range("E10").Value = "  +1.46%  "

# Row 11
$ws.Range("D11").Value = "'3.616.04"
$ws.Range("E11").Value = "  +6.07%  "

# Row 12
$ws.Range("D12").Value = "'42.46"
$ws.Range("E12").Value = "  -1.52%  "

# Row 13
$ws.Range("E13").Value = "  +0.86%  "

# Row 14
$ws.Range("D14").Value = "'6.29"
$ws.Range("E14").Value = "  +0.36%  "

# Row 15
$ws.Range("D15").Value = "'4.310.33"
$ws.Range("E15").Value = "  +6.20%  "

# Row 16
$ws.Range("D16").Value = "'95.630.35"
$ws.Range("E16").Value = "  +2.72%  "

# Row 17
$ws.Range("E17").Value = "  +3.10%  "

# Row 18
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'7.93"
$ws.Range("E18").Value = "  -3.68%  "

# Row 19
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "'3.621.43"
$ws.Range("E19").Value = "  +5.79%  "

# Row 20
$ws.Range("D20").Value = "'13.11"
$ws.Range("E20").Value = "  +12.63%  "

# Row 21
$ws.Range("D21").Value = "'18.02"
$ws.Range("E21").Value = "  +0.27%  "

# Row 22
$ws.Range("E22").Value = "  +5.79%  "

# Row 23
$ws.Range("D23").Value = "'0.482"
$ws.Range("E23").Value = "  -3.80%  "

# Row 24
$ws.Range("D24").Value = "'505.84"
$ws.Range("E24").Value = "  +1.80%  "

# Row 25
$ws.Range("E25").Value = "  +7.64%  "

# Row 26
$ws.Range("E26").Value = "  -2.07%  "

# Row 27
$ws.Range("D27").Value = "'95.77"
$ws.Range("E27").Value = "  +6.40%  "

# Row 28
$ws.Range("D28").Value = "'12.67"
$ws.Range("E28").Value = "  +5.53%  "

# Row 29
$ws.Range("D29").Value = "'3.815.89"
$ws.Range("E29").Value = "  +6.29%  "

# Row 30
$ws.Range("D30").Value = "'3.14"
$ws.Range("E30").Value = "  +15.91%  "

# Row 31
$ws.Range("D31").Value = "'11.32"
$ws.Range("E31").Value = "  +0.09%  "

# Row 32
$ws.Range("D32").Value = "'0.998"
$ws.Range("E32").Value = "  +0.03%  "

# Row 33
$ws.Range("E33").Value = "  +0.03%  "

# Row 34
$ws.Range("E34").Value = "  -1.11%  "

# Row 35
$ws.Range("B35").Value = "Cronos"
$ws.Range("C35").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D35").Value = "'0.177"
$ws.Range("E35").Value = "  +1.40%  "

# Row 36
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'31.94"
$ws.Range("E36").Value = "  +10.48%  "

# Row 37
$ws.Range("E37").Value = "  +1.72%  "

# Row 38
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'575.16"
$ws.Range("E38").Value = "  +2.80%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "'8.13"
$ws.Range("E39").Value = "  +9.15%  "

# Row 40
$ws.Range("D40").Value = "'1.48"
$ws.Range("E40").Value = "  +5.52%  "

# Row 41
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
$ws.Range("D42").Value = "'0.924"
$ws.Range("E42").Value = "  +3.21%  "

# Row 43
$ws.Range("E43").Value = "  -0.19%  "

# Row 44
$ws.Range("D44").Value = "'34.87"
$ws.Range("E44").Value = "  +44.30%  "

# Row 45
$ws.Range("D45").Value = "'1.72"
$ws.Range("E45").Value = "  +1.33%  "

# Row 46
$ws.Range("D46").Value = "'23.66"
$ws.Range("E46").Value = "  -0.09%  "

# Row 47
$ws.Range("D47").Value = "'5.65"
$ws.Range("E47").Value = "  +4.06%  "

# Row 48
$ws.Range("D48").Value = "'2.23"
$ws.Range("E48").Value = "  +6.12%  "

# Row 49
$ws.Range("D49").Value = "'0.0413"
$ws.Range("E49").Value = "  -0.18%  "

# Row 50
$ws.Range("D50").Value = "'3.51"
$ws.Range("E50").Value = "  -3.59%  "

# Row 51
$ws.Range("D51").Value = "'53.58"
$ws.Range("E51").Value = "  +1.25%  "
